$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new header value first
$ws.Range("H1").Value = "Save"

# Copy formatting from the neighboring header cell (G1) so the new
# "Save" header matches the style used by the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Re-assert the value in case paste affected it
$ws.Range("H1").Value = "Save"

# Add the new numeric value for the data row
$ws.Range("H2").Value = 1
